$wb = $excel.ActiveWorkbook

# --- Update conversion text on "Hoja1" sheet, cell A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 12.37 = 49611.71 pesos`n✅ 49611.71 pesos = 12.27 = 968.21 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update rate values on "tasas" sheet ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 80.86799999999999
$ws2.Range("O10").Value = 4012
$ws2.Range("N12").Value = 4042.99
$ws2.Range("O12").Value = 78.902
